$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TestcaseID values for the API test rows (new API methods / test cases added)
$ws.Range("A2").Value = "TC05"
$ws.Range("A3").Value = "TC03"
$ws.Range("A4").Value = "TC01"
$ws.Range("A5").Value = "TC02"

# Move the active selection to F2, matching the latest authored state of the sheet
$ws.Range("F2").Select() | Out-Null
